# Apply the automatic intraday output update:
#  - Update a handful of existing numeric/text cells on rows 2-5
#  - Append a brand new data row (row 6) for ticker TWG

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value (e.g. "-13.28%") into a cell without
# letting Excel's smart-parsing turn it into a numeric percentage.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
}

# --- Row 2 (CETX) updates -------------------------------------------------
$ws.Range("C2").Value = 164.82
$ws.Range("D2").Value = 4810000
Set-TextValue $ws.Range("G2") "-13.28%"

# --- Row 3 (FULC) updates -------------------------------------------------
$ws.Range("D3").Value = 703000000

# --- Row 4 (IBIO) updates -------------------------------------------------
$ws.Range("C4").Value = 60.83
Set-TextValue $ws.Range("G4") "-13.47%"

# --- Row 5 (PAVS) updates -------------------------------------------------
$ws.Range("C5").Value = 72.68
Set-TextValue $ws.Range("G5") "-22.63%"

# --- New row 6 (TWG) -------------------------------------------------------
$ws.Range("A6").Value = "TWG"
$ws.Range("B6").Value = 45999
$ws.Range("B6").NumberFormat = $ws.Range("B5").NumberFormat
$ws.Range("C6").Value = 58.24
$ws.Range("D6").Value = 13700000
$ws.Range("E6").Value = 410000
$ws.Range("F6").Value = 460000
Set-TextValue $ws.Range("G6") "131.44%"
$ws.Range("H6").Value = 34.02
$ws.Range("I6").Value = 0.02
$ws.Range("J6").Value = 3.84
$ws.Range("K6").Value = 9.550000000000001
$ws.Range("L6").Value = 9.51
$ws.Range("M6").Value = 26.36
$ws.Range("N6").Value = 8.619999999999999
$ws.Range("O6").Value = 19.35
$ws.Range("P6").Value = 13201340
$ws.Range("Q6").Value = "2025-12-08 14:43:00"
$ws.Range("R6").Value = "2025-12-08 09:49:00"
$ws.Range("S6").Value = 11.3
$ws.Range("T6").Value = 15.89
$ws.Range("U6").Value = 9.5
$ws.Range("V6").Value = 9.5
$ws.Range("W6").Value = 0
$ws.Range("X6").Value = "2025-12-08 04:20:00"
$ws.Range("Y6").Value = 9.76
$ws.Range("Z6").Value = 9.130000000000001
$ws.Range("AA6").Value = 3352010
$ws.Range("AB6").Value = 10.59
$ws.Range("AC6").Value = 9.130000000000001
$ws.Range("AD6").Value = 3494625
$ws.Range("AE6").Value = 10.99
$ws.Range("AF6").Value = 8.619999999999999
$ws.Range("AG6").Value = 3986094
$ws.Range("AH6").Value = 10.99
$ws.Range("AI6").Value = 8.619999999999999
$ws.Range("AJ6").Value = 4398310
$ws.Range("AK6").Value = 9.99
$ws.Range("AL6").Value = 11.55
$ws.Range("AM6").Value = 8.619999999999999
$ws.Range("AN6").Value = 4919318
$ws.Range("AO6").Value = 10.04
$ws.Range("AP6").Value = 11.55
$ws.Range("AQ6").Value = 8.619999999999999
$ws.Range("AR6").Value = 5084604
$ws.Range("AS6").Value = 11.21
$ws.Range("AT6").Value = 20
$ws.Range("AU6").Value = 8.619999999999999
$ws.Range("AV6").Value = 7936570
$ws.Range("AW6").Value = 19.05
